# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the cf7305e0-0aea-48fc-a6b8-86d8956c5559.md row across all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 7 corresponds to cf7305e0-0aea-48fc-a6b8-86d8956c5559.md
$wsOverview.Range("G7").Value = "2016-08-28 04:41:30"

# zh-cn sheet: row 7 corresponds to cf7305e0-0aea-48fc-a6b8-86d8956c5559.md
$wsZhCn.Range("H7").Value = "2016-08-28 04:41:26"

# de-de sheet: row 7 corresponds to cf7305e0-0aea-48fc-a6b8-86d8956c5559.md
$wsDeDe.Range("H7").Value = "2016-08-28 04:41:30"
